$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the Title values down by one row (B13->B14, B12->B13, B11->B12, B10->B11), then put
# the new entry's title into B10. Work bottom-up so we don't clobber values before copying.
$ws.Range("B14").Value = $ws.Range("B13").Value()
$ws.Range("B13").Value = $ws.Range("B12").Value()
$ws.Range("B12").Value = $ws.Range("B11").Value()
$ws.Range("B11").Value = $ws.Range("B10").Value()
$ws.Range("B10").Value = "Image classification: Cracks in concrete"

# New hyperlink cell for the inserted row (match the existing C2 hyperlink cell's formatting).
$ws.Range("C10").Value = "https://www.kaggle.com/code/vishnu0399/ensuring-structural-safety-crack-detection"
$ws.Hyperlinks.Add($ws.Range("C10"), "https://www.kaggle.com/code/vishnu0399/ensuring-structural-safety-crack-detection", "", "", "https://www.kaggle.com/code/vishnu0399/ensuring-structural-safety-crack-detection")
$ws.Range("C2").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row height for the newly-populated (wrapped) row.
$ws.Rows("10").RowHeight = 43.5

# Update the active selection.
$ws.Range("B5").Select()
